# Fix the mis-scaled / mis-aligned IFRS financial figures on the
# "company_list" sheet (한국석유공업): rows 2-6 get corrected per-column
# values, and the stale duplicate rows 7-9 have their data cells cleared
# (only the A/B/C label columns are kept).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value2 = 3793
$ws.Range("E2").Value2 = 104
$ws.Range("F2").Value2 = 104
$ws.Range("G2").Value2 = 95
$ws.Range("H2").Value2 = 67
$ws.Range("I2").Value2 = 67
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 2310
$ws.Range("L2").Value2 = 1138
$ws.Range("M2").Value2 = 1172
$ws.Range("N2").Value2 = 1172
$ws.Range("O2").Value2 = 1
$ws.Range("P2").Value2 = 33
$ws.Range("Q2").Value2 = 104
$ws.Range("R2").Value2 = -22
$ws.Range("S2").Value2 = -27
$ws.Range("T2").Value2 = 26
$ws.Range("U2").Value2 = 78
$ws.Range("V2").Value2 = 601
$ws.Range("W2").Value2 = 2.73
$ws.Range("X2").Value2 = 1.77
$ws.Range("Y2").Value2 = 5.85
$ws.Range("Z2").Value2 = 2.91
$ws.Range("AA2").Value2 = 97.06
$ws.Range("AB2").Value2 = 3445.17
$ws.Range("AC2").Value2 = 10191
$ws.Range("AD2").Value2 = 7.06
$ws.Range("AE2").Value2 = 179286
$ws.Range("AF2").Value2 = 0.4
$ws.Range("AG2").Value2 = 1000
$ws.Range("AH2").Value2 = 1.39
$ws.Range("AI2").Value2 = 9.79
$ws.Range("AJ2").Value2 = 655200

# Row 3 updates
$ws.Range("D3").Value2 = 3888
$ws.Range("E3").Value2 = 182
$ws.Range("F3").Value2 = 182
$ws.Range("G3").Value2 = 201
$ws.Range("H3").Value2 = 145
$ws.Range("I3").Value2 = 147
$ws.Range("J3").Value2 = -2
$ws.Range("K3").Value2 = 2402
$ws.Range("L3").Value2 = 1110
$ws.Range("M3").Value2 = 1292
$ws.Range("N3").Value2 = 1293
$ws.Range("O3").Value2 = -1
$ws.Range("P3").Value2 = 33
$ws.Range("Q3").Value2 = 132
$ws.Range("R3").Value2 = -30
$ws.Range("S3").Value2 = -93
$ws.Range("T3").Value2 = 34
$ws.Range("U3").Value2 = 98
$ws.Range("V3").Value2 = 520
$ws.Range("W3").Value2 = 4.67
$ws.Range("X3").Value2 = 3.73
$ws.Range("Y3").Value2 = 11.96
$ws.Range("Z3").Value2 = 6.16
$ws.Range("AA3").Value2 = 85.93
$ws.Range("AB3").Value2 = 3845.39
$ws.Range("AC3").Value2 = 22504
$ws.Range("AD3").Value2 = 4.44
$ws.Range("AE3").Value2 = 198072
$ws.Range("AF3").Value2 = 0.5
$ws.Range("AG3").Value2 = 1300
$ws.Range("AH3").Value2 = 1.3
$ws.Range("AI3").Value2 = 5.76
$ws.Range("AJ3").Value2 = 655200

# Row 4 updates
$ws.Range("D4").Value2 = 4027
$ws.Range("E4").Value2 = 169
$ws.Range("F4").Value2 = 169
$ws.Range("G4").Value2 = 172
$ws.Range("H4").Value2 = 118
$ws.Range("I4").Value2 = 120
$ws.Range("J4").Value2 = -2
$ws.Range("K4").Value2 = 2622
$ws.Range("L4").Value2 = 1232
$ws.Range("M4").Value2 = 1390
$ws.Range("N4").Value2 = 1393
$ws.Range("O4").Value2 = -3
$ws.Range("P4").Value2 = 33
$ws.Range("Q4").Value2 = 52
$ws.Range("R4").Value2 = -49
$ws.Range("S4").Value2 = -15
$ws.Range("T4").Value2 = 62
$ws.Range("U4").Value2 = -10
$ws.Range("V4").Value2 = 637
$ws.Range("W4").Value2 = 4.19
$ws.Range("X4").Value2 = 2.93
$ws.Range("Y4").Value2 = 8.92
$ws.Range("Z4").Value2 = 4.7
$ws.Range("AA4").Value2 = 88.68
$ws.Range("AB4").Value2 = 4176.96
$ws.Range("AC4").Value2 = 18277
$ws.Range("AD4").Value2 = 5.64
$ws.Range("AE4").Value2 = 215678
$ws.Range("AF4").Value2 = 0.48
$ws.Range("AG4").Value2 = 1500
$ws.Range("AH4").Value2 = 1.46
$ws.Range("AI4").Value2 = 8.09
$ws.Range("AJ4").Value2 = 655200

# Row 5 updates
$ws.Range("D5").Value2 = 4245
$ws.Range("E5").Value2 = 142
$ws.Range("F5").Value2 = 142
$ws.Range("G5").Value2 = 171
$ws.Range("H5").Value2 = 99
$ws.Range("I5").Value2 = 97
$ws.Range("J5").Value2 = 1
$ws.Range("K5").Value2 = 2848
$ws.Range("L5").Value2 = 1381
$ws.Range("M5").Value2 = 1467
$ws.Range("N5").Value2 = 1469
$ws.Range("O5").Value2 = -1
$ws.Range("P5").Value2 = 33
$ws.Range("Q5").Value2 = 62
$ws.Range("R5").Value2 = -118
$ws.Range("S5").Value2 = 78
$ws.Range("T5").Value2 = 20
$ws.Range("U5").Value2 = 41
$ws.Range("V5").Value2 = 728
$ws.Range("W5").Value2 = 3.35
$ws.Range("X5").Value2 = 2.33
$ws.Range("Y5").Value2 = 6.8
$ws.Range("Z5").Value2 = 3.61
$ws.Range("AA5").Value2 = 94.09
$ws.Range("AB5").Value2 = 4433.87
$ws.Range("AC5").Value2 = 14847
$ws.Range("AD5").Value2 = 6.5
$ws.Range("AE5").Value2 = 227761
$ws.Range("AF5").Value2 = 0.42
$ws.Range("AG5").Value2 = 2000
$ws.Range("AH5").Value2 = 2.07
$ws.Range("AI5").Value2 = 13.26
$ws.Range("AJ5").Value2 = 655200

# Row 6 updates
$ws.Range("D6").Value2 = 4813
$ws.Range("E6").Value2 = 83
$ws.Range("F6").Value2 = 83
$ws.Range("G6").Value2 = 77
$ws.Range("H6").Value2 = 43
$ws.Range("I6").Value2 = 44
$ws.Range("K6").Value2 = 3073
$ws.Range("L6").Value2 = 1596
$ws.Range("M6").Value2 = 1477
$ws.Range("N6").Value2 = 1479
$ws.Range("P6").Value2 = 33
$ws.Range("Q6").Value2 = -89
$ws.Range("R6").Value2 = -123
$ws.Range("S6").Value2 = 185
$ws.Range("T6").Value2 = 50
$ws.Range("U6").Value2 = -138
$ws.Range("V6").Value2 = 958
$ws.Range("W6").Value2 = 1.72
$ws.Range("X6").Value2 = 0.89
$ws.Range("Y6").Value2 = 2.99
$ws.Range("Z6").Value2 = 1.45
$ws.Range("AA6").Value2 = 108.04
$ws.Range("AB6").Value2 = 4519.12
$ws.Range("AC6").Value2 = 6717
$ws.Range("AD6").Value2 = 20.02
$ws.Range("AE6").Value2 = 236747
$ws.Range("AF6").Value2 = 0.57
$ws.Range("AG6").Value2 = 1500
$ws.Range("AH6").Value2 = 1.12
$ws.Range("AI6").Value2 = 21.3
$ws.Range("AJ6").Value2 = 655200

# Clear rows 7-9 data cells (D:AJ), keep A/B/C
$ws.Range("D7:AJ9").ClearContents()
